# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (fund-level holdings detail) right
#    before the "总计" (totals) summary sheet.
# 2) Insert a new top row into "总计" for the 2022-Q1 quarter and
#    renumber the existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q1" sheet, placed immediately before "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($total)
$ws.Name = "2022-Q1"

function Set-HeaderCell($cell, $text) {
    $cell.Value = $text
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.BorderAround(1)
}

function Set-IndexCell($cell, $n) {
    $cell.Value = $n
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.BorderAround(1)
}

Set-HeaderCell $ws.Range("B1") "基金代码"
Set-HeaderCell $ws.Range("C1") "基金名称"
Set-HeaderCell $ws.Range("D1") "基金规模"
Set-HeaderCell $ws.Range("E1") "股票总仓位"
Set-HeaderCell $ws.Range("F1") "仓位占比"
Set-HeaderCell $ws.Range("G1") "持有市值(亿元)"
Set-HeaderCell $ws.Range("H1") "仓位排名"

Set-IndexCell $ws.Range("A2") 0
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "001141"
$ws.Range("C2").Value = "泰达宏利创盈灵活配置混合 - A"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "3.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "27.59"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "5.40"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.1868"
$ws.Range("H2").Value = 1

Set-IndexCell $ws.Range("A3") 1
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "001142"
$ws.Range("C3").Value = "泰达宏利创盈灵活配置混合 - B"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "27.59"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "5.40"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0.0302"
$ws.Range("H3").Value = 1

# ---------------------------------------------------------------------
# Step 2: prepend the 2022-Q1 summary row to "总计" and renumber index
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows(2).Insert()
$tot.Range("B2:D2").ClearFormats()

Set-IndexCell $tot.Range("A2") 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 2
$tot.Range("D2").Value = 0.22

$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4
$tot.Range("A7").Value = 5

Write-Output "2022-Q1 data added"
